# Schritte.xlsx - "Sensoren auslesen, Zeitaufwand minimiert"
# Appends the two new log entries (2018-05-06) at the bottom of the
# documentation table on Tabelle1, continuing the existing Datum /
# Bilddateiname / Kurzbeschreibung layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 111 (A111:C111). Copy its date-cell formatting
# down onto the two new rows before filling them in, so the new date cells
# keep the same number format/style as the rest of column A.
$ws.Range("A111").Copy() | Out-Null
$ws.Range("A112").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A113").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 112: optimization timing results
$ws.Range("A112").Value = 43226
$ws.Range("B112").Value = "2018-05-06 1.jpg"
$ws.Range("C112").Value = "Die Optimierungen sind vorerst abgeschlossen. Eine Berechnungsdauer von ca. 6ms Konnte erzielt werden. Zwar keiner Erstrebten 4ms abeer immerhin besser als 125ms"

# Row 113: control-board revision note
$ws.Range("A113").Value = 43226
$ws.Range("B113").Value = "2018-05-06 2.pdf"
$ws.Range("C113").Value = "Die Optimierungen bedurften auch einigen Umstrukturierungen auf der Platine. Weshalb die 5. Revision der Kontrollschaltung gemalt wurde"

# Match the author's final selection state (C113 selected).
$ws.Range("C113").Select()
